$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-01 Tuesday", "2025-04-02 Wednesday"),
    @("237×6=1422", "508×9=4572"),
    @("650×9=5850", "202×9=1818"),
    @("710×4=2840", "885×6=5310"),
    @("877×2=1754", "630×4=2520"),
    @("794×4=3176", "863×4=3452"),
    @("304×9=2736", "832×8=6656"),
    @("651×8=5208", "278×2=556"),
    @("240×7=1680", "480×9=4320"),
    @("458×4=1832", "272×9=2448"),
    @("758×5=3790", "338×9=3042"),
    @("703×5=3515", "933×3=2799"),
    @("525×8=4200", "742×2=1484"),
    @("558×7=3906", "457×4=1828"),
    @("359×3=1077", "834×2=1668"),
    @("906×6=5436", "933×7=6531"),
    @("561×9=5049", "845×9=7605"),
    @("926×4=3704", "483×3=1449"),
    @("773×3=2319", "437×4=1748"),
    @("614×2=1228", "388×8=3104"),
    @("833×3=2499", "423×9=3807"),
    @("922×6=5532", "714×5=3570"),
    @("657×3=1971", "848×4=3392"),
    @("978×7=6846", "489×5=2445"),
    @("116×8=928", "102×3=306"),
    @("792×4=3168", "768×3=2304")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
